# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (property -> vocabulary term) - curated dimension terms
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("I2").Value = "iaest-measure:nacionalidad"

# Row 3 (medida / dim)
$ws.Range("D3").Value = "dim"
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"

# Row 4 (datatype / URI)
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"

# Row 5 (old mapping-file references) no longer needed - remove entirely
$ws.Rows(5).Delete()
